$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting psnr_nlm/psnr_gnlm/psnr_bm3d (A:C) to B:D.
# The native column-insert already carries over per-cell styles/widths for the
# shifted columns (and moves the stray styled I15/I49 marker cells to J15/J49).
$ws.Columns(1).Insert()

# New header for the inserted column + copy the header style (bold/centered,
# bordered) from the neighboring header cell so A1 matches B1:D1 formatting.
$ws.Range("A1").Value = "psnr_dual"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$arr = New-Object 'object[,]' 50,1
$arr[0,0] = 35.8737584
$arr[1,0] = 34.367576360000001
$arr[2,0] = 33.242258
$arr[3,0] = 33.343154589999997
$arr[4,0] = 33.563916460000002
$arr[5,0] = 34.215342800000002
$arr[6,0] = 34.116074310000002
$arr[7,0] = 31.042537549999999
$arr[8,0] = 33.56472806
$arr[9,0] = 33.0981399
$arr[10,0] = 35.24623072
$arr[11,0] = 31.517862399999999
$arr[12,0] = 32.850727190000001
$arr[13,0] = 34.761755149999999
$arr[14,0] = 32.956643939999999
$arr[15,0] = 33.438986409999998
$arr[16,0] = 32.105901719999999
$arr[17,0] = 31.231994050000001
$arr[18,0] = 30.78834406
$arr[19,0] = 32.322617319999999
$arr[20,0] = 33.302342230000001
$arr[21,0] = 33.033702609999999
$arr[22,0] = 33.590382599999998
$arr[23,0] = 36.23699688
$arr[24,0] = 33.95876268
$arr[25,0] = 35.835473909999997
$arr[26,0] = 33.485978719999999
$arr[27,0] = 35.654690799999997
$arr[28,0] = 34.308185530000003
$arr[29,0] = 32.26939883
$arr[30,0] = 32.40814125
$arr[31,0] = 31.899109549999999
$arr[32,0] = 31.695940619999998
$arr[33,0] = 32.164499650000003
$arr[34,0] = 32.139545040000002
$arr[35,0] = 33.182599750000001
$arr[36,0] = 32.395735899999998
$arr[37,0] = 31.84832566
$arr[38,0] = 34.264384319999998
$arr[39,0] = 35.393887710000001
$arr[40,0] = 34.114728999999997
$arr[41,0] = 33.84263035
$arr[42,0] = 32.577649950000001
$arr[43,0] = 34.566123670000003
$arr[44,0] = 31.49534079
$arr[45,0] = 32.836130439999998
$arr[46,0] = 32.961451619999998
$arr[47,0] = 32.63348088
$arr[48,0] = 34.758508130000003
$arr[49,0] = 33.247196799999998
$ws.Range("A2:A51").Value = $arr

$ws.Range("H12").Select()
